# Auto-generated Excel COM-interop script
# Applies the per-cell updates to the cryptocurrency price/volume listing
# on the active worksheet, matching the "Updated cryptos list" GitHub
# Actions commit (Wed Dec 13 17:47:06 UTC 2023).
#
# Every destination cell stores its value as literal text in the source
# workbook (t="inlineStr"), including price strings that look numeric
# (e.g. "251.55", "0.0935"). Excel's Range.Value setter auto-coerces such
# strings to real numbers, which would introduce float rounding noise
# (e.g. 251.55 -> 251.55000000000001) and change the cell type. To keep
# these cells as exact text -- matching the source -- each assignment
# forces the Text number format before the write and clears it again
# afterwards so no stray cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '42.039.04'
Set-TextValue 'E2' '  +2.49%  '
Set-TextValue 'D3' '2.210.93'
Set-TextValue 'E3' '  +1.69%  '
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '251.55'
Set-TextValue 'E5' '  +0.50%  '
Set-TextValue 'E6' '  +0.37%  '
Set-TextValue 'D7' '67.70'
Set-TextValue 'E7' '  +1.01%  '
Set-TextValue 'E8' '  -0.12%  '
Set-TextValue 'D9' '0.617'
Set-TextValue 'E9' '  +9.25%  '
Set-TextValue 'D10' '38.87'
Set-TextValue 'E10' '  +5.42%  '
Set-TextValue 'D11' '59.33'
Set-TextValue 'E11' '  +1.96%  '
Set-TextValue 'D12' '0.0935'
Set-TextValue 'E12' '  +0.63%  '
Set-TextValue 'E13' '  +0.71%  '
Set-TextValue 'E14' '  +0.18%  '
Set-TextValue 'D15' '2.544.03'
Set-TextValue 'E15' '  +1.87%  '
Set-TextValue 'D16' '0.868'
Set-TextValue 'E16' '  +1.10%  '
Set-TextValue 'D17' '14.46'
Set-TextValue 'E17' '  +0.57%  '
Set-TextValue 'D18' '2.210.42'
Set-TextValue 'E18' '  +1.59%  '
Set-TextValue 'D19' '41.951.38'
Set-TextValue 'E19' '  +2.57%  '
Set-TextValue 'E20' '  +1.90%  '
Set-TextValue 'D21' '72.24'
Set-TextValue 'E21' '  +1.02%  '
Set-TextValue 'E22' '  -0.62%  '
Set-TextValue 'D23' '230.96'
Set-TextValue 'E23' '  +0.21%  '
Set-TextValue 'D24' '2.01'
Set-TextValue 'E24' '  -2.35%  '
Set-TextValue 'E25' '  +1.24%  '
Set-TextValue 'D27' '11.13'
Set-TextValue 'E27' '  -4.14%  '
Set-TextValue 'D28' '2.40'
Set-TextValue 'E28' '  -2.42%  '
Set-TextValue 'E29' '  -1.04%  '
Set-TextValue 'D30' '2.26'
Set-TextValue 'E30' '  +4.42%  '
Set-TextValue 'D31' '166.92'
Set-TextValue 'E31' '  -1.42%  '
Set-TextValue 'E32' '  -0.33%  '
Set-TextValue 'B33' 'InternetComputer(DFINITY)'
Set-TextValue 'C33' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D33' '5.91'
Set-TextValue 'E33' '  +10.08%  '
Set-TextValue 'B34' 'Kaspa'
Set-TextValue 'C34' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D34' '0.121'
Set-TextValue 'E34' '  +3.88%  '
Set-TextValue 'D35' '0.0776'
Set-TextValue 'E35' '  +7.57%  '
Set-TextValue 'E36' '  +0.36%  '
Set-TextValue 'B37' 'InjectiveProtocol'
Set-TextValue 'C37' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D37' '26.04'
Set-TextValue 'E37' '  +2.50%  '
Set-TextValue 'B38' 'Filecoin'
Set-TextValue 'C38' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D38' '4.58'
Set-TextValue 'E38' '  +0.77%  '
Set-TextValue 'D39' '4.08'
Set-TextValue 'E39' '  +2.16%  '
Set-TextValue 'D40' '0.0313'
Set-TextValue 'E40' '  +6.01%  '
Set-TextValue 'E41' '  +0.92%  '
Set-TextValue 'B42' 'FTXToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D42' '5.18'
Set-TextValue 'E42' '  +7.90%  '
Set-TextValue 'B43' 'THORChain'
Set-TextValue 'C43' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D43' '5.64'
Set-TextValue 'E43' '  +0.39%  '
Set-TextValue 'B44' 'Celestia'
Set-TextValue 'C44' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D44' '11.96'
Set-TextValue 'E44' '  -1.07%  '
Set-TextValue 'D45' '61.15'
Set-TextValue 'E45' '  -4.62%  '
Set-TextValue 'E46' '  -2.49%  '
Set-TextValue 'E47' '  +0.06%  '
Set-TextValue 'D48' '0.0997'
Set-TextValue 'E48' '  -1.19%  '
Set-TextValue 'E49' '  -0.26%  '
Set-TextValue 'E50' '  +2.07%  '
